$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.020.57'
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").Value = '3.122.78'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.Value = "'594.59"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.86%  '
$c = $ws.Range("D6")
$c.Value = "'136.21"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -5.29%  '
$ws.Range("D8").Value = '3.117.82'
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("E10").Value = '  -3.91%  '
$c = $ws.Range("D11")
$c.Value = "'5.24"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.27%  '
$c = $ws.Range("D12")
$c.Value = "'0.455"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.64%  '
$c = $ws.Range("D13")
$c.Value = "'0.0000245"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -4.82%  '
$c = $ws.Range("D14")
$c.Value = "'34.13"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '3.637.18'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("E16").Value = '  +1.85%  '
$ws.Range("D17").Value = '63.036.84'
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").Value = '3.124.48'
$ws.Range("E18").Value = '  -0.95%  '
$ws.Range("E19").Value = '  -2.87%  '
$c = $ws.Range("D20")
$c.Value = "'472.35"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.65%  '
$c = $ws.Range("D21")
$c.Value = "'14.08"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.696"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("E23").Value = '  -0.29%  '
$c = $ws.Range("D24")
$c.Value = "'85.98"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.09%  '
$c = $ws.Range("D25")
$c.Value = "'12.85"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -4.59%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("B28").Value = 'NEARProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D28")
$c.Value = "'6.93"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.89%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D29")
$c.Value = "'7.88"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -6.37%  '
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("E33").Value = '  -6.28%  '
$ws.Range("E34").Value = '  -4.99%  '
$c = $ws.Range("D35")
$c.Value = "'1.07"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.83%  '
$ws.Range("E36").Value = '  -3.04%  '
$c = $ws.Range("D37")
$c.Value = "'51.93"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").Value = '0.0₃0697'
$ws.Range("E38").Value = '  -9.57%  '
$c = $ws.Range("D39")
$c.Value = "'0.0386"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.44%  '
$c = $ws.Range("D40")
$c.Value = "'418.36"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -6.48%  '
$c = $ws.Range("D41")
$c.Value = "'8.21"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").Value = '2.892.93'
$ws.Range("E42").Value = '  +0.99%  '
$c = $ws.Range("D43")
$c.Value = "'2.67"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -11.78%  '
$ws.Range("E44").Value = '  -5.90%  '
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("E46").Value = '  +0.03%  '
$c = $ws.Range("D47")
$c.Value = "'2.10"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -5.90%  '
$c = $ws.Range("D48")
$c.Value = "'25.41"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("E49").Value = '  -0.56%  '
$c = $ws.Range("D50")
$c.Value = "'2.24"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -7.10%  '
$c = $ws.Range("D51")
$c.Value = "'119.64"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.06%  '
